# Updated symbol list on Sat Dec 31 08:36:13 UTC 2022 with GitHub Actions
#
# Refreshes the crypto price/volume snapshot on Sheet1:
#  - numeric "Price" values (column D) are updated to the latest quotes
#  - two pairs of rows (11/20 and 41/42) had their coin ranking swapped,
#    so Coin/Link/Price/Volume columns are rewritten for those rows
#
# Price cells store numbers as plain text (e.g. "245.62", or "25.50" with
# a significant trailing zero) in the original workbook. Writing them with a
# leading apostrophe forces Excel to keep them as text instead of silently
# re-typing them as numbers; re-applying the Normal style afterwards drops
# the quote-prefix formatting flag so the cell style is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'25.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.143"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05592"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.489"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.027"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8175"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8471"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1339"
$ws.Range("D10").Style = "Normal"
# Row 11: MandalaExchangeToken -> LiechtensteinCryptoassetsExchange
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03234"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("D12").Value = "'0.02862"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.09389"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.001532"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.0005969"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.006140"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'3.532"
$ws.Range("D17").Style = "Normal"
# Row 20: LiechtensteinCryptoassetsExchange -> MandalaExchangeToken
$ws.Range("B20").Value = "MandalaExchangeToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D20").Value = "'0.06943"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19MandalaExchangeTokenMDX"
$ws.Range("D22").Value = "'3.747"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.04686"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Value = "'0.001246"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.004604"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.00009597"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.03654"
$ws.Range("D40").Style = "Normal"
# Row 41: BKEXToken -> KickToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.003383"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICK"
# Row 42: KickToken -> BKEXToken
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1356"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41BKEXTokenBKKBestin24h"
$ws.Range("D43").Value = "'0.002471"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.007776"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005322"
$ws.Range("D45").Style = "Normal"

Write-Host "Applied 40 cell updates to Sheet1"
